$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table ("Form_Responses") currently covers A1:K8 (header + 7 data rows).
# Add a new data row (row 9) with a new testimonial response.

$tbl = $ws.ListObjects.Item(1)

# Grow the table by one row; this also extends ref/autoFilter to K9.
$lr = $tbl.ListRows.Add()

# Copy the formatting (font, number format, alignment, etc.) of the last
# existing data row (row 8) down into the freshly added row 9 so the new
# row matches the look of the other data rows.
$ws.Range("A8:K8").Copy($ws.Range("A9:K9"))

# Fill in the new row's values.
$ws.Range("A9").Value = 45994.611641932876
$ws.Range("B9").Value = "aj.marsiglio@yahoo.com"
$ws.Range("C9").Value = "Before we met and you took over your role as supervisor, my position faced many challenges. The #1 challenge in my mind was the lack of a clear and defined processes for accomplishing my work. "
$ws.Range("D9").Value = "The challenge created chaos not just for my specific role but for my whole team. It was very frustrating coming to work not know if anything had changed since yesterday or if plans were going to change today. It made me reconsider if I wanted to stay in my role. "
$ws.Range("E9").Value = "I believe the most valuable piece of advice that you have given me is to advocate for myself and properly show what I have accomplished. Early on in my career I had my nose down, grinding away at work and projects but once goals were met and accomplished I simply moved on to the next without showing what I had accomplished. "
$ws.Range("F9").Value = "What sets you apart from others is your ability to try to look at things from multiple angles. From what I have seen you genuinely take the time to try to look at an issue or task from many angles and try to get a good representation of how everyone sees it, not just how you see it. I believe that has allowed you to have a better view and provide better, more honest guidance and assistance. "
$ws.Range("G9").Value = "I saw an increase in efficiency in my work, an increase of accreditation by showing my work, and an overall higher job satisfaction."
$ws.Range("H9").Value = "Working with you definitely helped my career trajectory as it has already impacted me as I am now a grade 7 and doing well. I believe I am setup to continue to move upwards as well. "
$ws.Range("I9").Value = "I would say he does not rush to judgement as he takes his time properly identifying the surface as well as underlying issues of any situation. He creates detailed plans and goals to attack those problems. Most importantly, he gives fair and honest feedback and advice on how to solve a problem or achieve a goal. "
$ws.Range("J9").Value = "Accelerate: getting a promotion or a leadership position"
$ws.Range("K9").Value = "Yes"

# Match the row height used by the other data rows.
$ws.Rows.Item(9).RowHeight = 22.5

# Keep the worksheet's hidden _FilterDatabase defined name in sync with the
# table/autofilter's new extent.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Form Responses 1'!`$A`$1:`$K`$9"
    }
}
